# Add a "date_time" column (C) to the export sheet, formatted as dd.mm.yyyy,
# filled with 2020-01-01 for every data row (matches serial date 43831).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the date number format to the whole new column before writing
# values/header so Excel doesn't invent an extra "general date" numFmt.
$ws.Range("C2:C4").NumberFormat = "dd.mm.yyyy"
$ws.Range("C1").NumberFormat = "dd.mm.yyyy"
$ws.Range("C1").HorizontalAlignment = -4152

$ws.Range("C1").Value = "date_time"
$ws.Range("C2").Value = "2020-01-01"
$ws.Range("C3").Value = "2020-01-01"
$ws.Range("C4").Value = "2020-01-01"

[void]$ws.Range("C4").Select()
